# Update column G ("K") values on Sheet1, rows 2-14 (row 7 unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 3
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 0
$ws.Range("G14").Value = 2
